$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text
$ws.Range("A1").Value = "number"

# Update existing phone number
$ws.Range("A2").Value = "01833184048"

# Add new phone number row, matching the text style used by the other number cells
$ws.Range("A4").Value = "01833184049"
$ws.Range("A4").NumberFormat = "@"

# Update active selection to the newly added cell
$ws.Range("A4").Select()
